$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: header row of the score matrix, referencing each agent's name
$ws.Range("C18").Formula = "=B19"
$ws.Range("D18").Formula = "=B20"
$ws.Range("E18").Formula = "=B21"
$ws.Range("F18").Formula = "=B22"
$ws.Range("G18").Formula = "=B23"
$ws.Range("H18").Formula = "=B24"
$ws.Range("I18").Formula = "=B25"
$ws.Range("J18").Formula = "=B26"
$ws.Range("K18").Formula = "=B27"
$ws.Range("L18").Formula = "=B28"
$ws.Range("M18").Formula = "=B29"
$ws.Range("N18").Formula = "=B30"
$ws.Range("O18").Formula = "=B31"
$ws.Range("P18").Formula = "=B32"
$ws.Range("Q18").Formula = "=B33"

# Row 19: Tit for Tat
$ws.Range("E19").Value = 600
$ws.Range("M19").Value = 283
$ws.Range("N19").Value = 283
$ws.Range("Q19").Value = 436

# Row 21: Nydegger
$ws.Range("C21").Value = 600
$ws.Range("M21").Value = 429
$ws.Range("N21").Value = 564
$ws.Range("Q21").Value = 339

# Row 29: Feld
$ws.Range("C29").Value = 288
$ws.Range("E29").Value = 714
$ws.Range("N29").Value = 286
$ws.Range("Q29").Value = 498

# Row 30: Joss
$ws.Range("C30").Value = 288
$ws.Range("E30").Value = 624
$ws.Range("M30").Value = 286
$ws.Range("Q30").Value = 460

# Row 33: Random
$ws.Range("C33").Value = 441
$ws.Range("E33").Value = 774
$ws.Range("M33").Value = 373
$ws.Range("N33").Value = 415

# Update the view: scrolled position and active selection
$win = $excel.ActiveWindow
$win.ScrollRow = 15
$win.ScrollColumn = 1
$ws.Range("F30").Select()
